$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value2 = 4924.0713
$ws.Range("I64").Value2 = 4536.5557
$ws.Range("K64").Value2 = 4536.5557
$ws.Range("M64").Value2 = -4288.5557
$ws.Range("H67").Value2 = 4924.0713
$ws.Range("I67").Value2 = 4536.5557
$ws.Range("K67").Value2 = 4536.5557
$ws.Range("M67").Value2 = -3678.5557
$ws.Range("H69").Value2 = 12999.833
$ws.Range("I69").Value2 = 6572
$ws.Range("J69").Value2 = 21998.8
$ws.Range("K69").Value2 = 19716
$ws.Range("L69").Value2 = 65996.39999999999
$ws.Range("M69").Value2 = -18842
$ws.Range("N69").Value2 = -67744.39999999999
$ws.Range("H70").Value2 = 942008.5600000001
$ws.Range("I70").Value2 = 2036070.1
$ws.Range("J70").Value2 = 4241.4287
$ws.Range("K70").Value2 = 6108210.300000001
$ws.Range("L70").Value2 = 12724.2861
$ws.Range("M70").Value2 = -6107940.300000001
$ws.Range("N70").Value2 = -13264.2861
$ws.Range("H72").Value2 = 12999.833
$ws.Range("I72").Value2 = 6572
$ws.Range("J72").Value2 = 21998.8
$ws.Range("K72").Value2 = 59148
$ws.Range("L72").Value2 = 197989.2
$ws.Range("M72").Value2 = -54780
$ws.Range("N72").Value2 = -206725.2
$ws.Range("H73").Value2 = 942008.5600000001
$ws.Range("I73").Value2 = 2036070.1
$ws.Range("J73").Value2 = 4241.4287
$ws.Range("K73").Value2 = 6108210.300000001
$ws.Range("L73").Value2 = 12724.2861
$ws.Range("M73").Value2 = -6107274.300000001
$ws.Range("N73").Value2 = -14596.2861
$ws.Range("H132").Value2 = 2725.3076
$ws.Range("I132").Value2 = 2196.3547
$ws.Range("K132").Value2 = 6589.0641
$ws.Range("M132").Value2 = -4059.0641
$ws.Range("H137").Value2 = 25002190
$ws.Range("I137").Value2 = 38463176
$ws.Range("J137").Value2 = 3214.1428
$ws.Range("K137").Value2 = 115389528
$ws.Range("L137").Value2 = 9642.428400000001
$ws.Range("M137").Value2 = -115386978
$ws.Range("N137").Value2 = -14742.4284
$ws.Range("H138").Value2 = 5060.121
$ws.Range("J138").Value2 = 7815.3335
$ws.Range("L138").Value2 = 23446.0005
$ws.Range("N138").Value2 = -33726.00049999999
$ws.Range("H141").Value2 = 27782136
$ws.Range("I141").Value2 = 29414418
$ws.Range("K141").Value2 = 88243254
$ws.Range("M141").Value2 = -88238074

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value2 = 500
$ws.Range("J10").Value2 = 500
$ws.Range("L10").Value2 = 500
$ws.Range("N10").Value2 = -840
$ws.Range("H11").Value2 = 866666.7
$ws.Range("I11").Value2 = 866666.7
$ws.Range("K11").Value2 = 866666.7
$ws.Range("M11").Value2 = -866522.7
$ws.Range("H61").Value2 = 21518042
$ws.Range("I61").Value2 = 23335980
$ws.Range("K61").Value2 = 23335980
$ws.Range("M61").Value2 = -23335768
$ws.Range("H133").Value2 = 74600.2
$ws.Range("J133").Value2 = 0
$ws.Range("L133").Value2 = 0
$ws.Range("N133").ClearContents() | Out-Null
$ws.Range("H136").Value2 = 21518042
$ws.Range("I136").Value2 = 23335980
$ws.Range("K136").Value2 = 70007940
$ws.Range("M136").Value2 = -70005390

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value2 = 318.33334
$ws.Range("I12").Value2 = 318.33334
$ws.Range("K12").Value2 = 318.33334
$ws.Range("M12").Value2 = -148.33334
$ws.Range("H31").Value2 = 16132448
$ws.Range("I31").Value2 = 24393292
$ws.Range("J31").Value2 = 4132.1904
$ws.Range("K31").Value2 = 24393292
$ws.Range("L31").Value2 = 4132.1904
$ws.Range("M31").Value2 = -24392997
$ws.Range("N31").Value2 = -4722.1904
$ws.Range("H34").Value2 = 16132448
$ws.Range("I34").Value2 = 24393292
$ws.Range("J34").Value2 = 4132.1904
$ws.Range("K34").Value2 = 24393292
$ws.Range("L34").Value2 = 4132.1904
$ws.Range("M34").Value2 = -24393090
$ws.Range("N34").Value2 = -4536.1904
$ws.Range("H58").Value2 = 2473.276
$ws.Range("I58").Value2 = 2067.92
$ws.Range("J58").Value2 = 5006.75
$ws.Range("K58").Value2 = 2067.92
$ws.Range("L58").Value2 = 5006.75
$ws.Range("M58").Value2 = -1864.92
$ws.Range("N58").Value2 = -5412.75
$ws.Range("H136").Value2 = 2473.276
$ws.Range("I136").Value2 = 2067.92
$ws.Range("J136").Value2 = 5006.75
$ws.Range("K136").Value2 = 6203.76
$ws.Range("L136").Value2 = 15020.25
$ws.Range("M136").Value2 = -3653.76
$ws.Range("N136").Value2 = -20120.25

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value2 = 25.4
$ws.Range("I11").Value2 = 25.4
$ws.Range("K11").Value2 = 25.4
$ws.Range("M11").Value2 = 113.6
$ws.Range("H14").Value2 = 4699.5
$ws.Range("I14").Value2 = 4699.5
$ws.Range("K14").Value2 = 4699.5
$ws.Range("M14").Value2 = -4531.5
$ws.Range("H122").Value2 = 3004790.8
$ws.Range("I122").Value2 = 3670971
$ws.Range("J122").Value2 = 6979
$ws.Range("K122").Value2 = 11012913
$ws.Range("L122").Value2 = 20937
$ws.Range("M122").Value2 = -11010463
$ws.Range("N122").Value2 = -25837

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value2 = 1917799
$ws.Range("I93").Value2 = 690.8421
$ws.Range("K93").Value2 = 690.8421
$ws.Range("M93").Value2 = 557.1579
$ws.Range("H132").Value2 = 2768.7795
$ws.Range("I132").Value2 = 1798.1063
$ws.Range("J132").Value2 = 6570.5835
$ws.Range("K132").Value2 = 5394.3189
$ws.Range("L132").Value2 = 19711.7505
$ws.Range("M132").Value2 = -2864.3189
$ws.Range("N132").Value2 = -24771.7505

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value2 = 5005500
$ws.Range("I10").Value2 = 5005500
$ws.Range("J10").Value2 = 0
$ws.Range("K10").Value2 = 5005500
$ws.Range("L10").Value2 = 0
$ws.Range("M10").Value2 = -5005331
$ws.Range("N10").ClearContents() | Out-Null
$ws.Range("H17").Value2 = 18199.8
$ws.Range("I17").Value2 = 18249.75
$ws.Range("K17").Value2 = 18249.75
$ws.Range("M17").Value2 = -18077.75
$ws.Range("H62").Value2 = 9144.333000000001
$ws.Range("J62").Value2 = 12259.8
$ws.Range("L62").Value2 = 12259.8
$ws.Range("N62").Value2 = -13507.8
$ws.Range("H65").Value2 = 9144.333000000001
$ws.Range("J65").Value2 = 12259.8
$ws.Range("L65").Value2 = 61299
$ws.Range("N65").Value2 = -67539
$ws.Range("H81").Value2 = 4130.5557
$ws.Range("I81").Value2 = 3397
$ws.Range("K81").Value2 = 6794
$ws.Range("M81").Value2 = -5733
$ws.Range("H84").Value2 = 4130.5557
$ws.Range("I84").Value2 = 3397
$ws.Range("K84").Value2 = 33970
$ws.Range("M84").Value2 = -28666
$ws.Range("H113").Value2 = 795.6
$ws.Range("I113").Value2 = 822.4286
$ws.Range("J113").Value2 = 733
$ws.Range("K113").Value2 = 2467.2858
$ws.Range("L113").Value2 = 2199
$ws.Range("M113").Value2 = -297.2857999999997
$ws.Range("N113").Value2 = -6539

